$d = $word.ActiveDocument

# Normal style: east-Asian font DejaVu Sans -> Tahoma
# (docDefaults/rPrDefault in styles.xml normally mirrors this in the
# source LibreOffice-authored template, but Word's object model has no
# member representing <w:docDefaults>/<w:rPrDefault> directly, so only
# the style-level change is reachable through COM automation.)
$d.Styles("Normal").Font.NameFarEast = "Tahoma"

# Heading style: east-Asian font DejaVu Sans -> Tahoma
$d.Styles("Heading").Font.NameFarEast = "Tahoma"

# List, Caption and Index styles previously had no explicit rFonts in
# their own rPr (they inherited everything). Setting the complex-script
# (bidi) font explicitly materializes a <w:rFonts w:cs="DejaVu Sans"/>
# element in each style's rPr, matching the target markup.
$d.Styles("List").Font.NameBi = "DejaVu Sans"
$d.Styles("Caption").Font.NameBi = "DejaVu Sans"
$d.Styles("Index").Font.NameBi = "DejaVu Sans"
